$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value first (so the new shared string "alice 333" is appended before the others)
$ws.Range("B8").Value = "alice 333"

# Add new header columns
$ws.Range("E1").Value = "folder_id"
$ws.Range("F1").Value = "ok2"

# Add new data values
$ws.Range("F3").Value = 2
$ws.Range("E5").Value = 222
$ws.Range("E7").Value = 444

# Update selection to match final state
$ws.Range("F8").Select()
